# Fill in the "Day" column (B) for each date row with the weekday name,
# based on the date value already present in column A (September 2022).
# September 1, 2022 was a Thursday, so the weekday cycles every 7 days
# starting from Thursday.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dayNames = @("Thursday", "Friday", "Saturday", "Sunday", "Monday", "Tuesday", "Wednesday")

# Dates live in column A, rows 6 through 35 (day-of-month 1..30 for September 2022).
for ($row = 6; $row -le 35; $row++) {
    $dateValue = $ws.Cells.Item($row, 1).Value2
    if ($dateValue -ne $null) {
        $idx = ([int]$dateValue - 1) % 7
        $ws.Cells.Item($row, 2).Value = $dayNames[$idx]
    }
}
